$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column O (year 2021) entirely - shifts nothing else, just removes
# the last data column from the table.
$ws.Columns("O").Delete()

# Correct the 2018 value (L5) from 1.6 to 1.7
$ws.Range("L5").Value = 1.7

# Correct the 2020 value (N5) from 3.1 to 1.6
$ws.Range("N5").Value = 1.6

# Update the active selection to reflect where the editor left off
$ws.Range("P6").Select()
